$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty placeholder rows 11-13 with real product
# ledger entries (date, product code, quantity, invoice number, staff).
# Row 11 keeps its pre-existing "date" style (s=1); the other columns use
# the sheet's default style, same as the existing data rows above.

$ws.Range("A11").Value = 45239
$ws.Range("B11").Value = "Apple"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 2
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = 2311090004
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "Ji"
$ws.Range("E11").Style = "Normal"

$ws.Range("A12").Value = 45239
$ws.Range("B12").Value = "Banana"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 3
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 2311090004
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "Ji"
$ws.Range("E12").Style = "Normal"

$ws.Range("A13").Value = 45239
$ws.Range("B13").Value = "Berry "
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 2
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = 2311090004
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "Ji"
$ws.Range("E13").Style = "Normal"
